# Finance_Data_Chemicals.xlsx - "Add files via upload" edit replay
#
# The commit swaps one row of the BSE Chemicals screen (Raghav Productivity
# Enhancers) for a different company (Agarwal Industrial Corporation),
# nudges a bunch of Market-Cap (column Q) values to a freshly-recomputed
# precision, makes the header row taller / wraps it, clears the sheet's
# AutoFilter + sort state, and moves the active selection to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 17 used to be "Raghav Productivity Enhancers Limited" (rank 41).
#    That row is removed; what used to be row 18 ("Punjab Chemicals and
#    Crop Protection Limited", rank 42) slides up into row 17 untouched.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = 42
$ws.Range("B17").Value = "Punjab Chemicals and Crop Protection Limited (BSE:506618)"
$ws.Range("C17").Value = 4985957
$ws.Range("D17").Value = "Diversified Chemicals"

# ---------------------------------------------------------------------
# 2) Row 18 becomes a brand-new entity: "Agarwal Industrial Corporation
#    Limited" (rank 43), with its own Entity ID and Industry
#    Classification. The rest of the row (ratios in E..P) is unchanged.
# ---------------------------------------------------------------------
$ws.Range("A18").Value = 43
$ws.Range("B18").Value = "Agarwal Industrial Corporation Limited (BSE:531921)"
$ws.Range("C18").Value = 4986658
$ws.Range("D18").Value = "Commodity Chemicals"

# ---------------------------------------------------------------------
# 3) Column Q (Market Cap) got re-pasted with a slightly different
#    floating point rounding on a couple dozen rows.
# ---------------------------------------------------------------------
$ws.Range("Q2").Value = 248695.89300000001
$ws.Range("Q3").Value = 248695.89300000001
$ws.Range("Q4").Value = 196607.054
$ws.Range("Q5").Value = 196607.054
$ws.Range("Q6").Value = 177240.82199999999
$ws.Range("Q7").Value = 109300.251
$ws.Range("Q8").Value = 54467.186099999999
$ws.Range("Q9").Value = 31058.971000000001
$ws.Range("Q11").Value = 23792.132099999999
$ws.Range("Q12").Value = 23792.132099999999
$ws.Range("Q13").Value = 22578.480100000001
$ws.Range("Q14").Value = 21858.3639
$ws.Range("Q15").Value = 20465.4244
$ws.Range("Q16").Value = 18590.850900000001
$ws.Range("Q17").Value = 18590.850900000001
$ws.Range("Q18").Value = 18590.850900000001
$ws.Range("Q22").Value = 10823.8951
$ws.Range("Q23").Value = 9894.1235199999992
$ws.Range("Q24").Value = 8540.09231
$ws.Range("Q25").Value = 7541.2983400000003
$ws.Range("Q26").Value = 5960.1290200000003
$ws.Range("Q28").Value = 5336.0935099999997
$ws.Range("Q31").Value = 2286.1864399999999
$ws.Range("Q33").Value = 1093.5493300000001
$ws.Range("Q34").Value = 1093.5493300000001
$ws.Range("Q35").Value = 453.76184899999998
$ws.Range("Q37").Value = 426.22183000000001

# ---------------------------------------------------------------------
# 4) Header row grows to fit wrapped text (row height -> 72pt).
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 72

# ---------------------------------------------------------------------
# 5) Drop the AutoFilter (and its saved sort state) that lived on A1:Q1.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false

# ---------------------------------------------------------------------
# 6) Move the live selection from A34:E34 to B1.
# ---------------------------------------------------------------------
[void]$ws.Range("B1").Select()
